$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2
$ws.Range('D2').Value = '67.149.69'
$ws.Range('E2').Value = '  +0.03%  '

# Row 3
$ws.Range('D3').Value = '2.486.97'
$ws.Range('E3').Value = '  -0.14%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
Set-TextValue $ws.Range('D5') '583.87'
$ws.Range('E5').Value = '  -0.14%  '

# Row 6
Set-TextValue $ws.Range('D6') '171.16'
$ws.Range('E6').Value = '  +3.24%  '

# Row 8
$ws.Range('E8').Value = '  -0.72%  '

# Row 9
$ws.Range('D9').Value = '2.486.43'
$ws.Range('E9').Value = '  -0.18%  '

# Row 10
$ws.Range('E10').Value = '  +0.89%  '

# Row 11
$ws.Range('E11').Value = '  +0.02%  '

# Row 12
$ws.Range('E12').Value = '  -0.48%  '

# Row 13
$ws.Range('E13').Value = '  -2.17%  '

# Row 14
$ws.Range('D14').Value = '2.960.39'
$ws.Range('E14').Value = '  +0.53%  '

# Row 15
Set-TextValue $ws.Range('D15') '25.31'
$ws.Range('E15').Value = '  -2.29%  '

# Row 16
$ws.Range('D16').Value = '67.039.93'
$ws.Range('E16').Value = '  -0.04%  '

# Row 17
$ws.Range('E17').Value = '  -1.65%  '

# Row 18
$ws.Range('D18').Value = '2.493.05'
$ws.Range('E18').Value = '  +0.35%  '

# Row 19
Set-TextValue $ws.Range('D19') '10.99'
$ws.Range('E19').Value = '  -5.23%  '

# Row 20
Set-TextValue $ws.Range('D20') '7.40'
$ws.Range('E20').Value = '  -5.67%  '

# Row 21
Set-TextValue $ws.Range('D21') '348.52'
$ws.Range('E21').Value = '  -3.03%  '

# Row 22
Set-TextValue $ws.Range('D22') '4.03'
$ws.Range('E22').Value = '  -1.99%  '

# Row 23
$ws.Range('E23').Value = '  -0.14%  '

# Row 24
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D24') '68.38'
$ws.Range('E24').Value = '  -3.09%  '

# Row 25
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D25') '4.23'
$ws.Range('E25').Value = '  -4.23%  '

# Row 26
Set-TextValue $ws.Range('D26') '1.79'
$ws.Range('E26').Value = '  -2.79%  '

# Row 27
$ws.Range('E27').Value = '  -1.69%  '

# Row 28
$ws.Range('E28').Value = '  +1.18%  '

# Row 29
$ws.Range('E29').Value = '  +0.03%  '

# Row 30
$ws.Range('E30').Value = '  -2.98%  '

# Row 31
Set-TextValue $ws.Range('D31') '510.81'
$ws.Range('E31').Value = '  +2.45%  '

# Row 32
Set-TextValue $ws.Range('D32') '7.72'
$ws.Range('E32').Value = '  -3.98%  '

# Row 33
Set-TextValue $ws.Range('D33') '1.23'
$ws.Range('E33').Value = '  -2.87%  '

# Row 34
$ws.Range('E34').Value = '  -3.78%  '

# Row 35
Set-TextValue $ws.Range('D35') '0.999'
$ws.Range('E35').Value = '  -0.05%  '

# Row 36
Set-TextValue $ws.Range('D36') '160.73'
$ws.Range('E36').Value = '  +0.59%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.118'
$ws.Range('E37').Value = '  -6.68%  '

# Row 38
$ws.Range('E38').Value = '  +0.72%  '

# Row 39
$ws.Range('E39').Value = '  -4.06%  '

# Row 40
$ws.Range('E40').Value = '  -5.32%  '

# Row 41
$ws.Range('E41').Value = '  -2.34%  '

# Row 42
$ws.Range('E42').Value = '  -0.05%  '

# Row 43
$ws.Range('E43').Value = '  -2.19%  '

# Row 44
$ws.Range('E44').Value = '  -2.62%  '

# Row 45
Set-TextValue $ws.Range('D45') '2.36'
$ws.Range('E45').Value = '  -4.06%  '

# Row 46
Set-TextValue $ws.Range('D46') '38.78'
$ws.Range('E46').Value = '  -1.37%  '

# Row 47
Set-TextValue $ws.Range('D47') '142.86'
$ws.Range('E47').Value = '  +0.67%  '

# Row 48
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D48') '3.45'
$ws.Range('E48').Value = '  -4.62%  '

# Row 49
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D49') '0.515'
$ws.Range('E49').Value = '  -4.29%  '

# Row 50
$ws.Range('E50').Value = '  -4.45%  '

# Row 51
Set-TextValue $ws.Range('D51') '0.0729'
$ws.Range('E51').Value = '  -1.23%  '
